$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "emp16.farida.m.abdelaziz@gmail.com"

for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 6).Value = $newValue
}
